# Generate Report for Archive
#
# Two files (002c40aa-5f5f-46ba-a936-b3c68b872ec8.md and
# 19c10cf3-27d7-4022-856e-bbda2dcb46f7.md) have moved from "Ready for
# handoff" to "In Translation". Update the Status on every sheet that
# tracks it: the per-language Status column on the "zh-cn" and "de-de"
# detail sheets, and the mirrored status columns on the "Overview" sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: File Name in column A, per-language status in B (zh-cn) and C (de-de) ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "In Translation"
$ov.Range("C3").Value = "In Translation"
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"

# --- zh-cn sheet: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

# --- de-de sheet: Status is column C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
